$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37.4534161490683
$ws.Range("C2").Value = 30.761421319797
$ws.Range("D2").Value = 33.8961038961039
$ws.Range("E2").Value = 33.4640522875817
